$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added at the top of this block (row 111),
# pushing the previously-existing rows 111-141 down to 112-142.
$ws.Rows.Item(111).Insert()

$ws.Cells.Item(111, 1).Value = 11
$ws.Cells.Item(111, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(111, 3).Value = 'Bíobío'
$ws.Cells.Item(111, 4).Value = 44754
$ws.Cells.Item(111, 5).Value = 8
$ws.Cells.Item(111, 6).Value = 100112043
$ws.Cells.Item(111, 7).Value = 'Pepino ensalada'
$ws.Cells.Item(111, 8).Value = 'Sin especificar'
$ws.Cells.Item(111, 9).Value = 'Primera'
$ws.Cells.Item(111, 10).Value = 100
$ws.Cells.Item(111, 11).Value = 19000
$ws.Cells.Item(111, 12).Value = 20000
$ws.Cells.Item(111, 13).Value = 19500
$ws.Cells.Item(111, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(111, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(111, 16).Value = 325
$ws.Cells.Item(111, 17).Value = 60
$ws.Cells.Item(111, 18).Value = 'Hortaliza'
